$wb = $excel.ActiveWorkbook

# --- Sheet "Simple Fields" ---
$ws1 = $wb.Worksheets.Item("Simple Fields")
$ws1.Range("A2").Value = ".)03nrie"
$ws1.Range("C2").Value = "34516 I =8 90 1"
$ws1.Range("D2").Value = "5th A ve."
$ws1.Range("E2").Value = "Brooklyn. NY"
$ws1.Range("H2").Value = "l 1 567 a 90 23 67 a 9 1 _"

# --- Sheet "Simple Fields - Formatted" ---
$ws2 = $wb.Worksheets.Item("Simple Fields - Formatted")
$ws2.Range("A2").Value = ".)03nrie"
$ws2.Range("C2").Value = "Key,Value`n`"Value`",`"`""
$ws2.Range("C2").WrapText = $true
$ws2.Range("D2").Value = "Key,Value`n`"Address Line 1`",`"5th A ve`""
$ws2.Range("E2").Value = "Brooklyn. NY"
$ws2.Range("H2").Value = "Key,Value`n`"Value`",`"`""
$ws2.Range("H2").WrapText = $true
$ws2.Rows.Item(2).AutoFit()

# --- Sheet "dependents" ---
$ws3 = $wb.Worksheets.Item("dependents")
$ws3.Range("A2").Value = "riot,lry Doe"
$ws3.Range("C3").Value = ". son"
$ws3.Range("D2").ClearContents()
$ws3.Range("D3").ClearContents()

# --- Sheet "dependents - Formatted" ---
$ws4 = $wb.Worksheets.Item("dependents - Formatted")
$ws4.Range("A2").Value = "Key,Value`n`"Given Name`",`"riot,lry`"`n`"Last Name`",`"Doe`""
$ws4.Range("C3").Value = ". son"
$ws4.Range("D2").ClearContents()
$ws4.Range("D3").ClearContents()
$ws4.Rows.Item(2).AutoFit()
